$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.233789666666667
$ws.Range("H2").Value = 3.701369
$ws.Range("I2").Value = 0.0001664233864291757
$ws.Range("J2").Value = 0.0001664233864291757
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1465046666666667
$ws.Range("N2").Value = 0.439514
$ws.Range("O2").Value = 0.07745172725947863
$ws.Range("P2").Value = 0.07745172725947864
$ws.Range("Q2").Value = 0.1807559438517778
$ws.Range("R2").Value = 1.626803494666
$ws.Range("S2").Value = 0.00001288977873531134
$ws.Range("T2").Value = 0.00001288977873531134

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.233789666666667
$ws.Range("H3").Value = 3.701369
$ws.Range("I3").Value = 0.0001664233864291757
$ws.Range("J3").Value = 0.0001664233864291757
$ws.Range("N3").Value = 4.707498
$ws.Range("O3").Value = 0.8295614045753745
$ws.Range("P3").Value = 0.8295614045753745
$ws.Range("Q3").Value = 1.936020796084667
$ws.Range("R3").Value = 17.424187164762
$ws.Range("S3").Value = 0.0001380584182003773
$ws.Range("T3").Value = 0.0001380584182003773

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.233789666666667
$ws.Range("H4").Value = 3.701369
$ws.Range("I4").Value = 0.0001664233864291757
$ws.Range("J4").Value = 0.0001664233864291757
$ws.Range("O4").Value = 0.09298686816514685
$ws.Range("P4").Value = 0.09298686816514684
$ws.Range("Q4").Value = 0.2170116757332222
$ws.Range("R4").Value = 1.953105081599
$ws.Range("S4").Value = 0.00001547518949348705
$ws.Range("T4").Value = 0.00001547518949348705

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("I5").Value = 0.9827534361704352
$ws.Range("J5").Value = 0.9827534361704352
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1465046666666667
$ws.Range("N5").Value = 0.439514
$ws.Range("O5").Value = 0.07745172725947863
$ws.Range("P5").Value = 0.07745172725947864
$ws.Range("Q5").Value = 1067.389197756542
$ws.Range("R5").Value = 9606.50277980888
$ws.Range("S5").Value = 0.07611595110158799
$ws.Range("T5").Value = 0.076115951101588

$ws.Range("I6").Value = 0.9827534361704352
$ws.Range("J6").Value = 0.9827534361704352
$ws.Range("N6").Value = 4.707498
$ws.Range("O6").Value = 0.8295614045753745
$ws.Range("P6").Value = 0.8295614045753745
$ws.Range("R6").Value = 102892.2687853965
$ws.Range("S6").Value = 0.815254320860822
$ws.Range("T6").Value = 0.815254320860822

$ws.Range("I7").Value = 0.9827534361704352
$ws.Range("J7").Value = 0.9827534361704352
$ws.Range("O7").Value = 0.09298686816514685
$ws.Range("P7").Value = 0.09298686816514684
$ws.Range("S7").Value = 0.09138316420802532
$ws.Range("T7").Value = 0.09138316420802531

$ws.Range("I8").Value = 0.01708014044313564
$ws.Range("J8").Value = 0.01708014044313564
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1465046666666667
$ws.Range("N8").Value = 0.439514
$ws.Range("O8").Value = 0.07745172725947863
$ws.Range("P8").Value = 0.07745172725947864
$ws.Range("Q8").Value = 18.55110013780289
$ws.Range("R8").Value = 166.959901240226
$ws.Range("S8").Value = 0.001322886379155332
$ws.Range("T8").Value = 0.001322886379155333

$ws.Range("I9").Value = 0.01708014044313564
$ws.Range("J9").Value = 0.01708014044313564
$ws.Range("N9").Value = 4.707498
$ws.Range("O9").Value = 0.8295614045753745
$ws.Range("P9").Value = 0.8295614045753745
$ws.Range("S9").Value = 0.01416902529635226
$ws.Range("T9").Value = 0.01416902529635226

$ws.Range("I10").Value = 0.01708014044313564
$ws.Range("J10").Value = 0.01708014044313564
$ws.Range("O10").Value = 0.09298686816514685
$ws.Range("P10").Value = 0.09298686816514684
$ws.Range("S10").Value = 0.001588228767628047
$ws.Range("T10").Value = 0.001588228767628047
